# Update handback timestamps on the "Generate Report for Handback" run.
# Column D = "Correspond Handoff Datetime", Column G = "Correspond Handback DateTime"
# Row 5 in both the "zh-cn" and "de-de" sheets corresponds to the
# 53a4fc39-b36d-457a-8cd3-e1e4dc9b40c8... handback entry.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-27 03:06:04"
$wsZhCn.Range("G5").Value = "2016-01-27 03:06:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-27 03:06:18"
$wsDeDe.Range("G5").Value = "2016-01-27 03:07:18"
